# Generate Report for Handoff
#
# Inserts two new localization-file rows (90bdb9d8-... and a73468fc-...)
# in between the existing "125b7fa1-..." row and the "fe8a6ab5-..." row
# on all three worksheets (Overview, zh-cn, de-de), pushing the
# "fe8a6ab5-..." row and the ".localization-config" row down by two rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item(1)

# Shift rows 3-4 ("fe8a6ab5-..." and ".localization-config") down to
# rows 5-6, leaving two blank rows (3 and 4) for the new entries.
$ws.Rows.Item(3).Insert()
$ws.Rows.Item(3).Insert()

# Row 3: 90bdb9d8-...
$ws.Range("A3").Value = "90bdb9d8-f14d-4e05-8ad1-65bf9d87a292.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"

# Row 4: a73468fc-...
$ws.Range("A4").Value = "a73468fc-c6a1-4fc6-ae9d-716cc413f21c.md"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = "Ready for handoff"

# Row 5 / Row 6 already carry the correct values from the row shift above
# (fe8a6ab5-... / .localization-config) - re-assert them explicitly so the
# content is correct regardless of how the shift landed.
$ws.Range("A5").Value = "fe8a6ab5-e15e-4df8-bf77-43d86abfd6b0.md"
$ws.Range("B5").Value = "Ready for handoff"
$ws.Range("C5").Value = "Ready for handoff"

$ws.Range("A6").Value = ".localization-config"
$ws.Range("B6").Value = "Not to be localized"
$ws.Range("C6").Value = "Not to be localized"

# Rebuild hyperlinks (row insert does not shift existing hyperlink refs).
$ws.Cells.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/7c4861cb7c2ccd396a73bb746e3ab7feb02ee8af/e2e/125b7fa1-ae8f-46ad-ab80-0ed7a0f66097.md", "", "", "125b7fa1-ae8f-46ad-ab80-0ed7a0f66097.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/d5fe6daee1aeeebafdd4118a77ee3cb2bbf5b059/e2e/90bdb9d8-f14d-4e05-8ad1-65bf9d87a292.md", "", "", "90bdb9d8-f14d-4e05-8ad1-65bf9d87a292.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/711b9945ecc75c721fd52433ca3f59cf799db265/e2e/a73468fc-c6a1-4fc6-ae9d-716cc413f21c.md", "", "", "a73468fc-c6a1-4fc6-ae9d-716cc413f21c.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/9593a224600e44cd8a54e2270c2e7157d0d5e2e4/e2e/fe8a6ab5-e15e-4df8-bf77-43d86abfd6b0.md", "", "", "fe8a6ab5-e15e-4df8-bf77-43d86abfd6b0.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/7c4861cb7c2ccd396a73bb746e3ab7feb02ee8af/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item(2)

$ws.Rows.Item(3).Insert()
$ws.Rows.Item(3).Insert()

# Row 3: 90bdb9d8-...
$ws.Range("A3").Value = "90bdb9d8-f14d-4e05-8ad1-65bf9d87a292.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "90bdb9d8-f14d-4e05-8ad1-65bf9d87a292.d5fe6daee1aeeebafdd4118a77ee3cb2bbf5b059.zh-cn.xlf"
$ws.Range("D3").Value = "2016-03-10 16:35:00"
$ws.Range("E3").Clear()
$ws.Range("F3").Clear()
$ws.Range("G3").Value = "0001-01-01 00:00:00"
$ws.Range("H3").Value = "Include"

# Row 4: a73468fc-...
$ws.Range("A4").Value = "a73468fc-c6a1-4fc6-ae9d-716cc413f21c.md"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = "a73468fc-c6a1-4fc6-ae9d-716cc413f21c.711b9945ecc75c721fd52433ca3f59cf799db265.zh-cn.xlf"
$ws.Range("D4").Value = "2016-03-10 16:35:00"
$ws.Range("E4").Clear()
$ws.Range("F4").Clear()
$ws.Range("G4").Value = "0001-01-01 00:00:00"
$ws.Range("H4").Value = "Include"

# Row 5: fe8a6ab5-... (shifted down from row 3)
$ws.Range("A5").Value = "fe8a6ab5-e15e-4df8-bf77-43d86abfd6b0.md"
$ws.Range("B5").Value = "Ready for handoff"
$ws.Range("C5").Value = "fe8a6ab5-e15e-4df8-bf77-43d86abfd6b0.b4c1648a74164d19502e406201443f14edf6a32d.zh-cn.xlf"
$ws.Range("D5").Value = "2016-03-10 16:33:39"
$ws.Range("G5").Value = "0001-01-01 00:00:00"
$ws.Range("H5").Value = "Include"

# Row 6: .localization-config (shifted down from row 4)
$ws.Range("A6").Value = ".localization-config"
$ws.Range("B6").Value = "Not to be localized"
$ws.Range("D6").Value = "0001-01-01 00:00:00"
$ws.Range("G6").Value = "0001-01-01 00:00:00"
$ws.Range("H6").Value = "Ignored"

$ws.Cells.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/7c4861cb7c2ccd396a73bb746e3ab7feb02ee8af/e2e/125b7fa1-ae8f-46ad-ab80-0ed7a0f66097.md", "", "", "125b7fa1-ae8f-46ad-ab80-0ed7a0f66097.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/11fc11f582a6164b74e3c018efbd82d5247f3f3e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/125b7fa1-ae8f-46ad-ab80-0ed7a0f66097.b512a4d960b68a30a3ab7a177a8ba77a6fe5a5db.zh-cn.xlf", "", "", "125b7fa1-ae8f-46ad-ab80-0ed7a0f66097.b512a4d960b68a30a3ab7a177a8ba77a6fe5a5db.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/3d8e364882d065cb3183e41a5cbe6b9d176f14df/e2e/125b7fa1-ae8f-46ad-ab80-0ed7a0f66097.md", "", "", "125b7fa1-ae8f-46ad-ab80-0ed7a0f66097.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/fb7c2b58ed18e4de0bce92b92b01fee263460ec1/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/125b7fa1-ae8f-46ad-ab80-0ed7a0f66097.b512a4d960b68a30a3ab7a177a8ba77a6fe5a5db.zh-cn.xlf", "", "", "125b7fa1-ae8f-46ad-ab80-0ed7a0f66097.b512a4d960b68a30a3ab7a177a8ba77a6fe5a5db.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/d5fe6daee1aeeebafdd4118a77ee3cb2bbf5b059/e2e/90bdb9d8-f14d-4e05-8ad1-65bf9d87a292.md", "", "", "90bdb9d8-f14d-4e05-8ad1-65bf9d87a292.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d5fe6daee1aeeebafdd4118a77ee3cb2bbf5b059/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/90bdb9d8-f14d-4e05-8ad1-65bf9d87a292.d5fe6daee1aeeebafdd4118a77ee3cb2bbf5b059.zh-cn.xlf", "", "", "90bdb9d8-f14d-4e05-8ad1-65bf9d87a292.d5fe6daee1aeeebafdd4118a77ee3cb2bbf5b059.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/711b9945ecc75c721fd52433ca3f59cf799db265/e2e/a73468fc-c6a1-4fc6-ae9d-716cc413f21c.md", "", "", "a73468fc-c6a1-4fc6-ae9d-716cc413f21c.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/711b9945ecc75c721fd52433ca3f59cf799db265/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a73468fc-c6a1-4fc6-ae9d-716cc413f21c.711b9945ecc75c721fd52433ca3f59cf799db265.zh-cn.xlf", "", "", "a73468fc-c6a1-4fc6-ae9d-716cc413f21c.711b9945ecc75c721fd52433ca3f59cf799db265.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/9593a224600e44cd8a54e2270c2e7157d0d5e2e4/e2e/fe8a6ab5-e15e-4df8-bf77-43d86abfd6b0.md", "", "", "fe8a6ab5-e15e-4df8-bf77-43d86abfd6b0.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0f4e281a1cdcedf95f896eb71e65c81f27f4f8d3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/fe8a6ab5-e15e-4df8-bf77-43d86abfd6b0.b4c1648a74164d19502e406201443f14edf6a32d.zh-cn.xlf", "", "", "fe8a6ab5-e15e-4df8-bf77-43d86abfd6b0.b4c1648a74164d19502e406201443f14edf6a32d.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/7c4861cb7c2ccd396a73bb746e3ab7feb02ee8af/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item(3)

$ws.Rows.Item(3).Insert()
$ws.Rows.Item(3).Insert()

# Row 3: 90bdb9d8-...
$ws.Range("A3").Value = "90bdb9d8-f14d-4e05-8ad1-65bf9d87a292.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "90bdb9d8-f14d-4e05-8ad1-65bf9d87a292.d5fe6daee1aeeebafdd4118a77ee3cb2bbf5b059.de-de.xlf"
$ws.Range("D3").Value = "2016-03-10 16:35:06"
$ws.Range("E3").Clear()
$ws.Range("F3").Clear()
$ws.Range("G3").Value = "0001-01-01 00:00:00"
$ws.Range("H3").Value = "Include"

# Row 4: a73468fc-...
$ws.Range("A4").Value = "a73468fc-c6a1-4fc6-ae9d-716cc413f21c.md"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = "a73468fc-c6a1-4fc6-ae9d-716cc413f21c.711b9945ecc75c721fd52433ca3f59cf799db265.de-de.xlf"
$ws.Range("D4").Value = "2016-03-10 16:35:06"
$ws.Range("E4").Clear()
$ws.Range("F4").Clear()
$ws.Range("G4").Value = "0001-01-01 00:00:00"
$ws.Range("H4").Value = "Include"

# Row 5: fe8a6ab5-... (shifted down from row 3)
$ws.Range("A5").Value = "fe8a6ab5-e15e-4df8-bf77-43d86abfd6b0.md"
$ws.Range("B5").Value = "Ready for handoff"
$ws.Range("C5").Value = "fe8a6ab5-e15e-4df8-bf77-43d86abfd6b0.b4c1648a74164d19502e406201443f14edf6a32d.de-de.xlf"
$ws.Range("D5").Value = "2016-03-10 16:33:43"
$ws.Range("G5").Value = "0001-01-01 00:00:00"
$ws.Range("H5").Value = "Include"

# Row 6: .localization-config (shifted down from row 4)
$ws.Range("A6").Value = ".localization-config"
$ws.Range("B6").Value = "Not to be localized"
$ws.Range("D6").Value = "0001-01-01 00:00:00"
$ws.Range("G6").Value = "0001-01-01 00:00:00"
$ws.Range("H6").Value = "Ignored"

$ws.Cells.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/7c4861cb7c2ccd396a73bb746e3ab7feb02ee8af/e2e/125b7fa1-ae8f-46ad-ab80-0ed7a0f66097.md", "", "", "125b7fa1-ae8f-46ad-ab80-0ed7a0f66097.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/468abc9da3edc8659891fe1397da2d14b2e02887/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/125b7fa1-ae8f-46ad-ab80-0ed7a0f66097.b512a4d960b68a30a3ab7a177a8ba77a6fe5a5db.de-de.xlf", "", "", "125b7fa1-ae8f-46ad-ab80-0ed7a0f66097.b512a4d960b68a30a3ab7a177a8ba77a6fe5a5db.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/2931a5aec41c97ef7b492ff9a62ee339251523fc/e2e/125b7fa1-ae8f-46ad-ab80-0ed7a0f66097.md", "", "", "125b7fa1-ae8f-46ad-ab80-0ed7a0f66097.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/5f485886866bdeb0ec193ee36b3cb2508d9b5765/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/125b7fa1-ae8f-46ad-ab80-0ed7a0f66097.b512a4d960b68a30a3ab7a177a8ba77a6fe5a5db.de-de.xlf", "", "", "125b7fa1-ae8f-46ad-ab80-0ed7a0f66097.b512a4d960b68a30a3ab7a177a8ba77a6fe5a5db.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/d5fe6daee1aeeebafdd4118a77ee3cb2bbf5b059/e2e/90bdb9d8-f14d-4e05-8ad1-65bf9d87a292.md", "", "", "90bdb9d8-f14d-4e05-8ad1-65bf9d87a292.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d5fe6daee1aeeebafdd4118a77ee3cb2bbf5b059/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/90bdb9d8-f14d-4e05-8ad1-65bf9d87a292.d5fe6daee1aeeebafdd4118a77ee3cb2bbf5b059.de-de.xlf", "", "", "90bdb9d8-f14d-4e05-8ad1-65bf9d87a292.d5fe6daee1aeeebafdd4118a77ee3cb2bbf5b059.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/711b9945ecc75c721fd52433ca3f59cf799db265/e2e/a73468fc-c6a1-4fc6-ae9d-716cc413f21c.md", "", "", "a73468fc-c6a1-4fc6-ae9d-716cc413f21c.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/711b9945ecc75c721fd52433ca3f59cf799db265/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a73468fc-c6a1-4fc6-ae9d-716cc413f21c.711b9945ecc75c721fd52433ca3f59cf799db265.de-de.xlf", "", "", "a73468fc-c6a1-4fc6-ae9d-716cc413f21c.711b9945ecc75c721fd52433ca3f59cf799db265.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/9593a224600e44cd8a54e2270c2e7157d0d5e2e4/e2e/fe8a6ab5-e15e-4df8-bf77-43d86abfd6b0.md", "", "", "fe8a6ab5-e15e-4df8-bf77-43d86abfd6b0.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cf500bfb72f2f0589048f92eeaf3aab47b6307b9/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/fe8a6ab5-e15e-4df8-bf77-43d86abfd6b0.b4c1648a74164d19502e406201443f14edf6a32d.de-de.xlf", "", "", "fe8a6ab5-e15e-4df8-bf77-43d86abfd6b0.b4c1648a74164d19502e406201443f14edf6a32d.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/7c4861cb7c2ccd396a73bb746e3ab7feb02ee8af/.localization-config", "", "", ".localization-config") | Out-Null

Write-Output "Generate Report for Handoff: rows inserted on all sheets"
